# "Add files via upload" - new 配音 (voice-over) column + copy-edits.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header for column E: 配音 (voice-over)
$ws.Range("E1").Value = "配音"

# Row 2 (林欣欣 line): add the voice-over file reference
$ws.Range("E2").Value = "GoWork\peiyin\1.wav"

# Row 3 used to be 顾诗琴; she's now voiced by 林欣欣 as well, reusing her
# icon, the dialogue line is tightened up, and a second voice-over file is
# added.
$ws.Range("A3").Value = "林欣欣"
$ws.Range("C3").Value = "同志，继续这样无所事事，我们是无法让人们获得解放的。"
$ws.Range("D3").Value = "GoWork\icons\111~1.ico"
$ws.Range("E3").Value = "GoWork\peiyin\2.wav"

# A couple of stray whitespace-only cells further down the sheet.
$ws.Range("F14").Value = "    "
$ws.Range("F16").Value = " "

# Leave the selection where the author left it.
$ws.Range("E3").Select()
